{"js": "// Update the course header:\n//   \"EG-247/EG-3068 Signals and Systems 2021-2022\"\n//   -> \"EG-247 Signals and Systems 2022-2023\"\n// (drop the second course code, roll the academic year forward by one)\n\nconst body = context.document.body;\n\n// 1) Drop the \"/EG-3068\" second course-code reference.\nconst codeMatches = body.search(\"/EG-3068\", { matchCase: true, matchWholeWord: false });\ncodeMatches.load(\"items\");\nawait context.sync();\n\nif (codeMatches.items.length > 0) {\n  codeMatches.items[0].delete();\n  await context.sync();\n}\n\n// 2) Roll the academic year forward: \"2021-2022\" -> \"2022-2023\".\nconst yearMatches = body.search(\"2021-2022\", { matchCase: true, matchWholeWord: false });\nyearMatches.load(\"items\");\nawait context.sync();\n\nif (yearMatches.items.length > 0) {\n  yearMatches.items[0].insertText(\"2022-2023\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Update the course header:\n#   \"EG-247/EG-3068 Signals and Systems 2021-2022\"\n#   -> \"EG-247 Signals and Systems 2022-2023\"\n# (drop the second course code, roll the academic year forward by one)\n\n$d = $word.ActiveDocument\n\n# 1) Drop the \"/EG-3068\" second course-code reference.\n$dropCode = $d.Content.Find\n$foundCode = $dropCode.Execute(\"/EG-3068\", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 2)\n\n# 2) Roll the academic year forward: \"2021-2022\" -> \"2022-2023\".\n$bumpYear = $d.Content.Find\n$foundYear = $bumpYear.Execute(\"2021-2022\", $false, $false, $false, $false, $false, $true, 1, $false, \"2022-2023\", 2)\n"}
